$d = $word.ActiveDocument

# Locate the sentence fragment "...aplicação na tela." and turn it into
# "...aplicação no ecrã." The original text lives in a single run; the
# target keeps the text "...aplicação n" in that run and moves the new
# "o ecrã." into a second, identically-formatted run - exactly as Word
# does when a user overwrites a selection in the middle of a run.

# Step 1: do the actual text substitution ("a tela." -> "o ecrã.").
$rng = $d.Content
$rng.Find.Execute("a tela.", $true, $false, $false, $false, $false, $true, 1, $false, "o ecrã.", 2)

# Step 2: force Word to split the run right before "o ecrã." by touching
# the run-level formatting of just the replaced tail. Re-find it, nudge a
# property to a different value (forcing the split into its own run), then
# restore it to the original value so the visible formatting is unchanged.
$tail = $d.Content
$tail.Find.Execute("o ecrã.")
$tail.Font.Size = 40
$tail.Font.Size = 8
